# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-28) used to be sorted in
# descending order (2008 .. 1908). It is now sorted ascending
# (1908 .. 2008), and the "Valor Mora" amounts in column F follow the
# same rows, so the 57600 value (originally attached to period 2008)
# now travels with period 2008 at its new position (row 28), while the
# other rows keep the 72000 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008")
$values  = @(72000,72000,72000,72000,72000,72000,72000,72000,72000,72000,72000,72000,57600)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $values[$i]
}
